# [fix] id into excel
# Corrects id_authorization values (column A) on Sheet1 and restores the
# workbook's last-saved view state (active cell / scroll position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix id_authorization (column A) values ---------------------------------
$ws.Range("A15").Value = 1612
$ws.Range("A31").Value = 1628
$ws.Range("A32").Value = 1628
$ws.Range("A46").Value = 1643
$ws.Range("A59").Value = 1656
$ws.Range("A60").Value = 1656
$ws.Range("A61").Value = 1656
$ws.Range("A66").Value = 1663
$ws.Range("A67").Value = 1663
$ws.Range("A69").Value = 1666
$ws.Range("A73").Value = 1670
$ws.Range("A76").Value = 1673

# --- Restore view state (scroll position + active cell/selection) ----------
$excel.Goto($ws.Range("D73"), $true) | Out-Null
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D73").Select() | Out-Null
